# Updated cryptos list - refresh price/volume figures and re-rank two swapped pairs.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '31.028.06'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +1.23%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.957.48'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -0.21%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9992'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.27'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.75%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9993'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.05%  '

$ws.Range('E7').Value = '  +0.80%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2961'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.43%  '

$ws.Range('E9').Value = '  +0.46%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.18'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.13%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '107.08'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -2.84%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.957.88'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.37%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.07817'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.99%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.468'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.00%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.7038'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.23%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '284.39'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.55%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '31.056.65'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.23%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.22'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.41%  '

$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007684'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.17%  '

$ws.Range('B20').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C20').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.220.29'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.82%  '

$ws.Range('E21').Value = '  -0.02%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.512'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -2.10%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.9993'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.02%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.501'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.61%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.800'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.30%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '170.40'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.35%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.99'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.67%  '

$ws.Range('E28').Value = '  +0.82%  '

$ws.Range('E29').Value = '  -0.45%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.402'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.35%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.583'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.71%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.600'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.41%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.460'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.28%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04933'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.83%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7611'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.41%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.173'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.30%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.726'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.17%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02014'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.14%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.698'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.52%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.566'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.32%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '77.31'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +9.84%  '

$ws.Range('E42').Value = '  +0.80%  '

$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4472'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.30%  '

$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.8860'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.04%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '109.39'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.29%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.117'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +9.53%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.9997'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.05%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '997.31'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +9.12%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.1261'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.86%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.329'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.23%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '35.90'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.60%  '
